# Regenerate the "K" column (G) with freshly computed strikeout-per-outing
# values (s_vals) instead of the previous Strike# figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values, one per data row (rows 2-56), in row order.
$s_vals = @(2,0,0,0,1,2,2,1,2,2,0,1,1,3,1,3,1,1,0,0,1,1,1,2,3,2,1,0,1,0,2,1,1,3,1,2,1,2,1,0,2,2,1,2,2,1,3,2,1,1,1,0,1,2,3)

$startRow = 2
for ($i = 0; $i -lt $s_vals.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $s_vals[$i]
}
